$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.570.69"
$ws.Range("E2").Value = "  +0.09%  "
$ws.Range("D3").Value = "1.753.48"
$ws.Range("E3").Value = "  +0.00%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.002"
$ws.Range("E4").Value = "  -0.21%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "324.36"
$ws.Range("E5").Value = "  +0.02%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.001"
$ws.Range("E6").Value = "  -0.13%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4497"
$ws.Range("E7").Value = "  +0.76%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3571"
$ws.Range("E8").Value = "  -1.16%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07470"
$ws.Range("E9").Value = "  -0.59%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "41.42"
$ws.Range("E10").Value = "  -1.87%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.084"
$ws.Range("E11").Value = "  -1.60%  "
$ws.Range("E12").Value = "  -0.05%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "20.77"
$ws.Range("E13").Value = "  +0.31%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.990"
$ws.Range("E14").Value = "  -0.87%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.155"
$ws.Range("E15").Value = "  -0.75%  "
$ws.Range("D16").Value = "1.755.95"
$ws.Range("E16").Value = "  +0.00%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "94.47"
$ws.Range("E17").Value = "  +1.85%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001056"
$ws.Range("E18").Value = "  -0.71%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06392"
$ws.Range("E19").Value = "  -0.42%  "
$ws.Range("E20").Value = "  -0.01%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "17.08"
$ws.Range("E21").Value = "  +0.15%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.737"
$ws.Range("E22").Value = "  -2.00%  "
$ws.Range("D23").Value = "27.611.49"
$ws.Range("E23").Value = "  +0.07%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.18"
$ws.Range("E24").Value = "  -0.48%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.083"
$ws.Range("E25").Value = "  -1.19%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "165.53"
$ws.Range("E26").Value = "  +1.95%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "20.16"
$ws.Range("E27").Value = "  -1.28%  "
$ws.Range("D28").Value = "1.953.72"
$ws.Range("E28").Value = "  -0.18%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.116"
$ws.Range("E29").Value = "  -0.49%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "125.56"
$ws.Range("E30").Value = "  -0.24%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.087"
$ws.Range("E31").Value = "  -0.11%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.09163"
$ws.Range("E32").Value = "  +1.07%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.654"
$ws.Range("E33").Value = "  +0.55%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.508"
$ws.Range("E34").Value = "  -0.67%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.02283"
$ws.Range("E35").Value = "  -0.68%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "11.70"
$ws.Range("E36").Value = "  -3.03%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.2091"
$ws.Range("E37").Value = "  -0.57%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.05997"
$ws.Range("E38").Value = "  +0.21%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.6274"
$ws.Range("E39").Value = "  -1.79%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "4.917"
$ws.Range("E40").Value = "  -0.28%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.178"
$ws.Range("E41").Value = "  -1.28%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.391"
$ws.Range("E42").Value = "  -0.15%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "7.789"
$ws.Range("E43").Value = "  -0.33%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "13.16"
$ws.Range("E44").Value = "  -0.95%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.717"
$ws.Range("E45").Value = "  +0.18%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.5854"
$ws.Range("E46").Value = "  -0.53%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "121.81"
$ws.Range("E47").Value = "  +0.13%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.932"
$ws.Range("E48").Value = "  -1.78%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.06886"
$ws.Range("E49").Value = "  +0.16%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.130"
$ws.Range("E50").Value = "  -2.36%  "
$ws.Range("E51").Value = "  -0.31%  "
